$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell value while preserving it as literal text
# (avoids Excel auto-converting numeric-looking strings to numbers/percentages)
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "278.11"
Set-TextValue $ws.Range("E2") "0.92%"
Set-TextValue $ws.Range("D3") "27.41"
Set-TextValue $ws.Range("E3") "0.28%"
Set-TextValue $ws.Range("D4") "4.841"
Set-TextValue $ws.Range("E4") "0.81%"
Set-TextValue $ws.Range("D5") "0.06366"
Set-TextValue $ws.Range("E5") "0.40%"
Set-TextValue $ws.Range("D6") "7.033"
Set-TextValue $ws.Range("E6") "1.07%"
Set-TextValue $ws.Range("D7") "1.287"
Set-TextValue $ws.Range("E7") "-5.74%"
Set-TextValue $ws.Range("D8") "0.8926"
Set-TextValue $ws.Range("E8") "1.73%"
Set-TextValue $ws.Range("D9") "0.1516"
Set-TextValue $ws.Range("E9") "-0.35%"
Set-TextValue $ws.Range("D10") "0.05844"
Set-TextValue $ws.Range("E10") "15.08%"
Set-TextValue $ws.Range("D11") "0.07453"
Set-TextValue $ws.Range("E11") "-0.76%"
Set-TextValue $ws.Range("D12") "0.02915"
Set-TextValue $ws.Range("E12") "-1.69%"
Set-TextValue $ws.Range("D13") "0.08972"
Set-TextValue $ws.Range("D14") "0.001595"
Set-TextValue $ws.Range("E14") "1.34%"
Set-TextValue $ws.Range("D15") "0.0006402"
Set-TextValue $ws.Range("E15") "-0.24%"
Set-TextValue $ws.Range("D16") "0.006176"
Set-TextValue $ws.Range("E16") "6.99%"
Set-TextValue $ws.Range("D17") "3.470"
Set-TextValue $ws.Range("E17") "0.64%"
Set-TextValue $ws.Range("D18") "3.303"
Set-TextValue $ws.Range("E18") "0.12%"
Set-TextValue $ws.Range("D19") "2.250"
Set-TextValue $ws.Range("E19") "-0.98%"
Set-TextValue $ws.Range("D21") "0.1349"
Set-TextValue $ws.Range("E21") "-0.29%"
Set-TextValue $ws.Range("D22") "3.896"
Set-TextValue $ws.Range("E22") "-0.10%"
$ws.Range("B23").Value = "CoinExToken"
$ws.Range("C23").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
Set-TextValue $ws.Range("D23") "0.04409"
Set-TextValue $ws.Range("E23") "0.34%"
$ws.Range("B24").Value = "ZBToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
Set-TextValue $ws.Range("D24") "0.1505"
Set-TextValue $ws.Range("E24") "9.05%"
Set-TextValue $ws.Range("D25") "0.001175"
Set-TextValue $ws.Range("E25") "0.16%"
Set-TextValue $ws.Range("D26") "0.004268"
Set-TextValue $ws.Range("E26") "10.32%"
Set-TextValue $ws.Range("D28") "0.0001179"
Set-TextValue $ws.Range("E28") "-1.61%"
Set-TextValue $ws.Range("E29") "-14.49%"
Set-TextValue $ws.Range("D40") "0.04022"
Set-TextValue $ws.Range("E40") "-4.39%"
Set-TextValue $ws.Range("D41") "0.006704"
Set-TextValue $ws.Range("E41") "-1.64%"
Set-TextValue $ws.Range("D42") "0.1412"
Set-TextValue $ws.Range("E42") "19.72%"
Set-TextValue $ws.Range("D43") "0.002058"
Set-TextValue $ws.Range("E43") "2.04%"
Set-TextValue $ws.Range("D44") "0.01117"
Set-TextValue $ws.Range("E44") "-2.73%"
Set-TextValue $ws.Range("D45") "0.00005541"
Set-TextValue $ws.Range("D46") "1.561"
Set-TextValue $ws.Range("E46") "5.01%"
Set-TextValue $ws.Range("E47") "-19.47%"
